# Update the "取得日時" (retrieved datetime) timestamps in the ランサーズ sheet
# for rows 2-6 (column A) from "2026-01-21 06:33:53" to "2026-01-21 06:42:09".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-21 06:42:09"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
